# Add two new keyword rows to the bottom of the sheet, matching the
# formatting pattern used by the other "highlighted" keyword rows
# (row height 18 / style index 1, font "Sakkal Majalla").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 41 and 42.
$ws.Range("A41").Value = "بوئینگ"
$ws.Range("A42").Value = "مناطق آزاد"

# Copy the cell formatting (style) from the last existing "highlighted"
# row (A40) onto the new rows so they reuse the same cell style (s="1")
# instead of creating new font/style table entries.
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A40").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the taller row height (18) used by the other highlighted rows.
$ws.Range("A41").EntireRow.RowHeight = 18
$ws.Range("A42").EntireRow.RowHeight = 18

# Move the active selection to the next empty row, as Excel does after
# the last data entry.
[void]$ws.Range("A43").Select()
